$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
  # Row 17
  $ws.Range("H17").Value = 898.3333
  $ws.Range("J17").Value = 898.3333
  $ws.Range("L17").Value = 2694.9999
  $ws.Range("N17").Value = -3030.9999
  # Row 43
  $ws.Range("H43").Value = 1458.6666
  $ws.Range("I43").Value = 1438.5
  $ws.Range("J43").Value = 1499
  $ws.Range("K43").Value = 1438.5
  $ws.Range("L43").Value = 1499
  $ws.Range("M43").Value = -1369.5
  $ws.Range("N43").Value = -1637
  # Row 53
  $ws.Range("H53").Value = 449.2
  $ws.Range("I53").Value = 199
  $ws.Range("K53").Value = 199
  $ws.Range("M53").Value = 438
  # Row 76
  $ws.Range("H76").Value = 969
  $ws.Range("J76").Value = 969
  $ws.Range("L76").Value = 969
  $ws.Range("N76").Value = -1599
  # Row 79
  $ws.Range("H79").Value = 969
  $ws.Range("J79").Value = 969
  $ws.Range("L79").Value = 969
  $ws.Range("N79").Value = -3153
  # Row 86
  $ws.Range("H86").Value = 191640
  $ws.Range("I86").Value = 1102.5
  $ws.Range("K86").Value = 1102.5
  $ws.Range("M86").Value = 20.5
  # Row 89
  $ws.Range("H89").Value = 191640
  $ws.Range("I89").Value = 1102.5
  $ws.Range("K89").Value = 5512.5
  $ws.Range("M89").Value = 103.5
  # Row 92
  $ws.Range("H92").Value = 1313.091
  $ws.Range("J92").Value = 1462.5
  $ws.Range("L92").Value = 1462.5
  $ws.Range("N92").Value = -3958.5
  # Row 111
  $ws.Range("H111").Value = 5471.875
  $ws.Range("I111").Value = 4890.4614
  $ws.Range("K111").Value = 14671.3842
  $ws.Range("M111").Value = -11604.3842
  # Row 125
  $ws.Range("H125").Value = 2033.3334
  $ws.Range("I125").Value = 550.5
  $ws.Range("K125").Value = 4954.5
  $ws.Range("M125").Value = -2494.5
  # Row 127
  $ws.Range("H127").Value = 2615.625
  $ws.Range("I127").Value = 2615.625
  $ws.Range("K127").Value = 7846.875
  $ws.Range("M127").Value = -2886.875
  # Row 132
  $ws.Range("H132").Value = 2260
  $ws.Range("I132").Value = 2395.6667
  $ws.Range("K132").Value = 7187.000100000001
  $ws.Range("M132").Value = -4657.000100000001
  # Row 140
  $ws.Range("H140").Value = 55000
  $ws.Range("J140").Value = 55000
  $ws.Range("L140").Value = 55000
  $ws.Range("N140").Value = -65360

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
  # Row 2
  $ws.Range("H2").Value = 1978.4736
  $ws.Range("I2").Value = 1385.5
  $ws.Range("K2").Value = 1385.5
  $ws.Range("M2").Value = -1272.5
  # Row 3
  $ws.Range("H3").Value = 5470
  $ws.Range("I3").Value = 2649.5
  $ws.Range("J3").Value = 11111
  $ws.Range("K3").Value = 2649.5
  $ws.Range("L3").Value = 11111
  $ws.Range("M3").Value = -2534.5
  $ws.Range("N3").Value = -11341
  # Row 116
  $ws.Range("H116").Value = 1978.4736
  $ws.Range("I116").Value = 1385.5
  $ws.Range("K116").Value = 1385.5
  $ws.Range("M116").Value = 908.5
  # Row 119
  $ws.Range("H119").Value = 39799
  $ws.Range("J119").Value = 39799
  $ws.Range("L119").Value = 39799
  $ws.Range("N119").Value = -49475

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
  # Row 3
  $ws.Range("H3").Value = 1978.4736
  $ws.Range("I3").Value = 1385.5
  $ws.Range("K3").Value = 1385.5
  $ws.Range("M3").Value = -1271.5
  # Row 87
  $ws.Range("H87").Value = 70321
  $ws.Range("I87").Value = 70321
  $ws.Range("K87").Value = 70321
  $ws.Range("M87").Value = -69073
  # Row 90
  $ws.Range("H90").Value = 70321
  $ws.Range("I90").Value = 70321
  $ws.Range("K90").Value = 210963
  $ws.Range("M90").Value = -204723
  # Row 118
  $ws.Range("H118").Value = 200000
  $ws.Range("J118").Value = 200000
  $ws.Range("L118").Value = 200000
  $ws.Range("N118").Value = -203314

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
  # Row 4
  $ws.Range("H4").Value = 500000
  $ws.Range("J4").Value = 500000
  $ws.Range("L4").Value = 500000
  $ws.Range("N4").Value = -500224
  # Row 31
  $ws.Range("H31").Value = 2154.5833
  $ws.Range("I31").Value = 2154.5833
  $ws.Range("K31").Value = 2154.5833
  $ws.Range("M31").Value = -1859.5833
  # Row 34
  $ws.Range("H34").Value = 2154.5833
  $ws.Range("I34").Value = 2154.5833
  $ws.Range("K34").Value = 2154.5833
  $ws.Range("M34").Value = -1952.5833

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
  # Row 81
  $ws.Range("H81").Value = 2400
  $ws.Range("J81").Value = 2400
  $ws.Range("L81").Value = 7200
  $ws.Range("N81").Value = -9446
  # Row 84
  $ws.Range("H84").Value = 2400
  $ws.Range("J84").Value = 2400
  $ws.Range("L84").Value = 21600
  $ws.Range("N84").Value = -32832
  # Row 107
  $ws.Range("H107").Value = 150
  $ws.Range("J107").Value = 150
  $ws.Range("L107").Value = 450
  $ws.Range("N107").Value = -4290

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
  # Row 80
  $ws.Range("H80").Value = 19049.166
  $ws.Range("J80").Value = 100000
  $ws.Range("L80").Value = 100000
  $ws.Range("N80").Value = -101996
  # Row 83
  $ws.Range("H83").Value = 19049.166
  $ws.Range("J83").Value = 100000
  $ws.Range("L83").Value = 500000
  $ws.Range("N83").Value = -509984

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
  # Row 7
  $ws.Range("H7").Value = 6396.9165
  $ws.Range("I7").Value = 2933.3333
  $ws.Range("K7").Value = 2933.3333
  $ws.Range("M7").Value = -2821.3333
  # Row 16
  $ws.Range("H16").Value = 1879.8
  $ws.Range("I16").Value = 1879.8
  $ws.Range("K16").Value = 1879.8
  $ws.Range("M16").Value = -1709.8
  # Row 35
  $ws.Range("H35").Value = 2000
  $ws.Range("I35").Value = 2000
  $ws.Range("J35").Value = 0
  $ws.Range("K35").Value = 2000
  $ws.Range("L35").Value = 0
  $ws.Range("M35").Value = -1664
  $ws.Range("N35").ClearContents()
  # Row 40
  $ws.Range("H40").Value = 3637.889
  $ws.Range("I40").Value = 3360.4614
  $ws.Range("K40").Value = 3360.4614
  $ws.Range("M40").Value = -3224.4614
  # Row 46
  $ws.Range("H46").Value = 3251
  $ws.Range("I46").Value = 3293.5
  $ws.Range("K46").Value = 3293.5
  $ws.Range("M46").Value = -3105.5
  # Row 126
  $ws.Range("H126").Value = 6396.9165
  $ws.Range("I126").Value = 2933.3333
  $ws.Range("K126").Value = 8799.999899999999
  $ws.Range("M126").Value = -6329.999899999999

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
  # Row 45
  $ws.Range("H45").Value = 16011.667
  $ws.Range("J45").Value = 16011.667
  $ws.Range("L45").Value = 16011.667
  $ws.Range("N45").Value = -16993.667
  # Row 112
  $ws.Range("H112").Value = 25387
  $ws.Range("J112").Value = 25387
  $ws.Range("L112").Value = 25387
  $ws.Range("N112").Value = -28341
  # Row 126
  $ws.Range("H126").Value = 1263.2727
  $ws.Range("I126").Value = 899.7778
  $ws.Range("K126").Value = 2699.3334
  $ws.Range("M126").Value = -229.3334
